$d = $word.ActiveDocument

# Locate the three consecutive "employee" user-story paragraphs:
#   13: "(10 points): As an employee, I want to be able to confirm that I have completed a pickup."
#   14: "(10 points): As an employee, I want all confirmed pickups to have a charge applied to the customer."
#   15: "(10 points): As an employee, I want to be able to select a customer profile and see their address
#        with a pin on a map (Google Maps API)."
#
# The edit moves the "_GoBack" bookmark from the start of paragraph 13 to the
# start of paragraph 15, and applies the same yellow highlight used by the
# other user-story bullets to paragraph 15's run (it previously had none).

$target = "(10 points): As an employee, I want to be able to select a customer profile and see their address with a pin on a map (Google Maps API)."

$p15 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $p15 = $p
        break
    }
}

# Remove the bookmark from its current location (start of the "confirm pickup" paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Re-create it at the start of the "select a customer profile" paragraph.
$newBmRange = $d.Range($p15.Range.Start, $p15.Range.Start)
$d.Bookmarks.Add("_GoBack", $newBmRange)

# Highlight the "select a customer profile" run in yellow, matching its siblings.
$p15.Range.HighlightColorIndex = 7
